# "Compilation post corrections France et Pierre"
#
# The "Directeur technique adjoint" role went from one person (with a
# placeholder phone number) to two co-holders, and the Cafeteria / Race
# Headquarters contacts (France Galarneau & Pierre Galarneau / Emmanuel
# Gilbert) received their real phone numbers instead of "TBD".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO")

# "Directeur technique adjoint" (singular) -> "Directeurs techniques
# adjoints" (plural) - role_fr column for the Assistant technical
# director(s) row.
$ws.Range("B5").Value = "Directeurs techniques adjoints"

# Cafeteria contacts (France Galarneau / Pierre Galarneau) now have phone
# numbers instead of "TBD".
$ws.Range("D9").Value = "\(819) 732-4038<br/>(819) 727-8510"

# Race Headquarters contact (Emmanuel Gilbert) now has a phone number
# instead of "TBD".
$ws.Range("D10").Value = "\(819) 860-6928"

# Selection left on B15 when the file was saved.
$ws.Activate()
$ws.Range("B15").Select()
